$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 7694481.5
$ws.Range("J17").Value = 8335530
$ws.Range("L17").Value = 25006590
$ws.Range("N17").Value = -25006926
# Row 76
$ws.Range("H76").Value = 2418689.5
$ws.Range("I76").Value = 3370.5881
$ws.Range("K76").Value = 3370.5881
$ws.Range("M76").Value = -3055.5881
# Row 79
$ws.Range("H79").Value = 2418689.5
$ws.Range("I79").Value = 3370.5881
$ws.Range("K79").Value = 3370.5881
$ws.Range("M79").Value = -2278.5881
# Row 80
$ws.Range("H80").Value = 12825996
$ws.Range("I80").Value = 2634
$ws.Range("J80").Value = 15230376
$ws.Range("K80").Value = 7902
$ws.Range("L80").Value = 45691128
$ws.Range("M80").Value = -6904
$ws.Range("N80").Value = -45693124
# Row 83
$ws.Range("H83").Value = 12825996
$ws.Range("I83").Value = 2634
$ws.Range("J83").Value = 15230376
$ws.Range("K83").Value = 23706
$ws.Range("L83").Value = 137073384
$ws.Range("M83").Value = -18714
$ws.Range("N83").Value = -137083368
# Row 107
$ws.Range("H107").Value = 994.4211
$ws.Range("I107").Value = 938.3077
$ws.Range("J107").Value = 1116
$ws.Range("K107").Value = 938.3077
$ws.Range("L107").Value = 1116
$ws.Range("M107").Value = 981.6923
$ws.Range("N107").Value = -4956
# Row 116
$ws.Range("H116").Value = 4564.1816
$ws.Range("J116").Value = 4967.3335
$ws.Range("L116").Value = 4967.3335
$ws.Range("N116").Value = -11851.3335
# Row 129
$ws.Range("H129").Value = 334029.56
$ws.Range("I129").Value = 317.4
$ws.Range("J129").Value = 400772
$ws.Range("K129").Value = 952.1999999999999
$ws.Range("L129").Value = 1202316
$ws.Range("M129").Value = 4047.8
$ws.Range("N129").Value = -1212316
# Row 134
$ws.Range("H134").Value = 47954.5
$ws.Range("J134").Value = 47954.5
$ws.Range("L134").Value = 47954.5
$ws.Range("N134").Value = -58094.5
# Row 138
$ws.Range("H138").Value = 2478.9722
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 2478.9722
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 7436.9166
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -17716.9166

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5696.9463
$ws.Range("I32").Value = 4409.0137
$ws.Range("K32").Value = 4409.0137
$ws.Range("M32").Value = -4122.0137
# Row 97
$ws.Range("H97").Value = 55556860
$ws.Range("I97").Value = 939.75
$ws.Range("K97").Value = 939.75
$ws.Range("M97").Value = -443.75
# Row 112
$ws.Range("H112").Value = 32019
$ws.Range("J112").Value = 32019
$ws.Range("L112").Value = 32019
$ws.Range("N112").Value = -34973

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 76
$ws.Range("H76").Value = 16000
$ws.Range("J76").Value = 16000
$ws.Range("L76").Value = 16000
$ws.Range("N76").Value = -16630
# Row 79
$ws.Range("H79").Value = 16000
$ws.Range("J79").Value = 16000
$ws.Range("L79").Value = 16000
$ws.Range("N79").Value = -18184

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3483.862
$ws.Range("I31").Value = 842.3125
$ws.Range("J31").Value = 6735
$ws.Range("K31").Value = 842.3125
$ws.Range("L31").Value = 6735
$ws.Range("M31").Value = -547.3125
$ws.Range("N31").Value = -7325
# Row 34
$ws.Range("H34").Value = 3483.862
$ws.Range("I34").Value = 842.3125
$ws.Range("J34").Value = 6735
$ws.Range("K34").Value = 842.3125
$ws.Range("L34").Value = 6735
$ws.Range("M34").Value = -640.3125
$ws.Range("N34").Value = -7139
# Row 58
$ws.Range("H58").Value = 20130.852
$ws.Range("I58").Value = 1687.0667
$ws.Range("J58").Value = 43185.582
$ws.Range("K58").Value = 1687.0667
$ws.Range("L58").Value = 43185.582
$ws.Range("M58").Value = -1484.0667
$ws.Range("N58").Value = -43591.582
# Row 86
$ws.Range("H86").Value = 24884.875
$ws.Range("J86").Value = 29813
$ws.Range("L86").Value = 29813
$ws.Range("N86").Value = -32059
# Row 89
$ws.Range("H89").Value = 24884.875
$ws.Range("J89").Value = 29813
$ws.Range("L89").Value = 149065
$ws.Range("N89").Value = -160297
# Row 132
$ws.Range("H132").Value = 3366.5881
$ws.Range("I132").Value = 1756.3636
$ws.Range("J132").Value = 6318.6665
$ws.Range("K132").Value = 5269.0908
$ws.Range("L132").Value = 18955.9995
$ws.Range("M132").Value = -2739.0908
$ws.Range("N132").Value = -24015.9995
# Row 134
$ws.Range("H134").Value = 1908.8889
$ws.Range("I134").Value = 1863.3334
$ws.Range("K134").Value = 5590.0002
$ws.Range("M134").Value = -3055.0002
# Row 136
$ws.Range("H136").Value = 20130.852
$ws.Range("I136").Value = 1687.0667
$ws.Range("J136").Value = 43185.582
$ws.Range("K136").Value = 5061.2001
$ws.Range("L136").Value = 129556.746
$ws.Range("M136").Value = -2511.2001
$ws.Range("N136").Value = -134656.746

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 852.3182
$ws.Range("I113").Value = 680.3
$ws.Range("J113").Value = 995.6667
$ws.Range("K113").Value = 2040.9
$ws.Range("L113").Value = 2987.0001
$ws.Range("M113").Value = 129.1000000000001
$ws.Range("N113").Value = -7327.0001
# Row 131
$ws.Range("H131").Value = 692.95
$ws.Range("I131").Value = 402.5
$ws.Range("J131").Value = 718.20654
$ws.Range("K131").Value = 1207.5
$ws.Range("L131").Value = 2154.61962
$ws.Range("M131").Value = 3832.5
$ws.Range("N131").Value = -12234.61962
# Row 138
$ws.Range("H138").Value = 2688.1538
$ws.Range("I138").Value = 2011.875
$ws.Range("J138").Value = 3770.2
$ws.Range("K138").Value = 6035.625
$ws.Range("L138").Value = 11310.6
$ws.Range("M138").Value = -895.625
$ws.Range("N138").Value = -21590.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4475150
$ws.Range("I70").Value = 4683.3335
$ws.Range("K70").Value = 4683.3335
$ws.Range("M70").Value = -4413.3335
# Row 73
$ws.Range("H73").Value = 4475150
$ws.Range("I73").Value = 4683.3335
$ws.Range("K73").Value = 4683.3335
$ws.Range("M73").Value = -3747.3335
# Row 113
$ws.Range("H113").Value = 5382.778
$ws.Range("I113").Value = 6439.263
$ws.Range("J113").Value = 2873.625
$ws.Range("K113").Value = 6439.263
$ws.Range("L113").Value = 2873.625
$ws.Range("M113").Value = -4269.263
$ws.Range("N113").Value = -7213.625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2030.9445
$ws.Range("I22").Value = 2574.5386
$ws.Range("J22").Value = 617.6
$ws.Range("K22").Value = 2574.5386
$ws.Range("L22").Value = 617.6
$ws.Range("M22").Value = -2279.5386
$ws.Range("N22").Value = -1207.6
# Row 27
$ws.Range("H27").Value = 2030.9445
$ws.Range("I27").Value = 2574.5386
$ws.Range("J27").Value = 617.6
$ws.Range("K27").Value = 2574.5386
$ws.Range("L27").Value = 617.6
$ws.Range("M27").Value = -2467.5386
$ws.Range("N27").Value = -831.6

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 9057.571
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 9057.571
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 9057.571
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -9283.571
